$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the worksheet relationship id is not controllable via the object model;
# focus on the data content changes described in the diff.

# Final account numbers, names and balances (in the new row order) after the edit.
$accounts = @("002606448", "005255637", "004204344", "004368468", "004313254", "004213139", "004329030", "004392159", "005696595", "004363260", "004855960", "005003629", "004220849", "004432579", "004508516", "004355790", "005040864", "003301389", "004207184", "005701765", "001368670", "004239387", "004431591", "004383268", "004384258", "004466350", "004536602", "003115072", "004462930", "004212132", "004809902", "004517080", "003249855", "004261201", "004424761", "004384167", "004754920", "008069255", "000827730", "005142611", "001719494", "005880251", "004563252", "004477812", "004877741", "002823185", "005173958", "004452507", "004212409", "004381194", "008032257", "004457389", "004749928", "004242237", "004027477", "004335251", "005924958", "003836362", "005092207", "004374891", "005068961", "005141215", "005890232", "005558076", "004691225", "004384131", "004382374", "004974089", "004207278", "004335144", "004334062", "004546050", "005079458", "004218542", "004552021", "004504449", "004643880", "004400640", "004426743", "008032413", "004332207", "005076418", "004693308", "005152037", "004508159", "001731007", "005981575", "004278212", "004381095", "000806386", "004332103", "004425965", "004290978", "004971448", "005245032", "004238164", "004752519", "004752615", "005902068", "004272426", "004413537", "005077648", "004806286", "005266369", "004240292", "005009992", "005591536", "004435987", "004211911", "004207374", "005055239", "004385806", "005927101", "004377415", "004230529", "004813134", "005305965", "004459461", "005018038", "004404724", "005616259", "005044389", "004350197", "004472760", "005715733", "004756968", "004243043", "004277637", "004398174", "008070544", "008090243", "004371857", "004388077", "004214604", "004467884", "005143579", "004204255", "004214592", "004920447", "008032597", "008037529", "004994036", "003497496", "004975924", "002894447", "005293480", "004547722", "004340984", "004260002", "004422594", "004455356", "004570632", "005878792", "004454365", "005685089", "005009922", "004264780", "004374943", "004551472", "004472076", "004335031", "005374916", "004216298", "008013889", "002828327", "004751154", "003895497", "005324981", "004308815", "004505474", "004530494", "004752494", "008007759", "004907688", "004228456", "004448501", "004912314", "004381415", "004754056", "004805269", "005268516", "004377713", "005135281", "000834301", "004539779", "008012870", "004165515", "008004995", "004328934", "004181486", "005170415", "005140667", "005022526", "000431814", "004360430", "004486497", "004115403", "005660155", "004223502", "004587511", "004380948", "004473942", "002687737", "003894173", "004453302", "004638738", "004278033", "005662526", "004432455", "002694089", "004357159", "004320840", "001000288", "005530256", "004451996", "005047946", "004223226", "004281300", "005274028", "004329229", "004213373", "004339183", "004870976", "000938440", "002878817", "004400000", "004976625", "005171652")
$names = @("MARCUS", "PATRICIA", "CLINEO", "AHMAD", "GUSTAVO", "LEONARDO", "DANIELA", "RODRIGO", "CLUBE", "LARISSA", "CLERIA", "ANDRE", "DULCE", "ANA", "EDUARDO", "MINEIA", "ANDRE", "EDMUNDO", "CRISTINA", "F", "THIAGO", "LUIZ", "MARIO", "LAURA", "PAULA", "RAQUEL", "TATIANY", "VICTOR", "WALTER", "JOAO", "PEDRO", "TATIANA", "MARINA", "ANA", "PEDRO", "DOUGLAS", "LUIS", "ANGELA", "LUCIANA", "GUILHERME", "LUIS", "LUIZ", "FERNANDO", "DIEGO", "LUIZ", "SIMONE", "VENIA", "DANIELA", "RAFAEL", "ALINNE", "SARA", "RAFAEL", "NADY", "MARIAH", "GABRIELA", "EDMUNDO", "TIAGO", "ISABELLA", "BRUNO", "RODRIGO", "JORGE", "KARINA", "TAYLA", "ALEXANDRE", "ANNA", "ANDRE", "THEOMAR", "CELIA", "CESAR", "EDMUNDO", "MERG", "LUIS", "JONAS", "JOSE", "MARIA", "KELMA", "GABRIEL", "FELIPE", "GABRIELLE", "VICTOR", "IRACY", "LEONARDO", "LAURA", "RODRIGO", "FELIPE", "GUILHERME", "GLAUCIANE", "LEONARDO", "THIAGO", "FERNANDA", "JOSE", "CAROLLINA", "LARISSA", "CLOVIS", "ROSA", "DANIELA", "MARCUS", "LUZIMAR", "VIVIANE", "RODRIGO", "CLAUDIA", "DUNAS", "VERA", "EG", "MARCO", "ALINE", "GUSTAVO", "MARCO", "ZENILDA", "ANGELICA", "NORMAN", "ANILSON", "SIMONE", "ANGELA", "LAIS", "MONICA", "SIDMAR", "INTERLAGOS", "ELAINE", "LEANDRO", "MARIA", "CLAUDIA", "GISELA", "SANDRA", "ADRIANA", "DANIELY", "SUELI", "LARA", "DANIELE", "MARINA", "GABRIEL", "NAZARETH", "WLADMIR", "MARIA", "ANA", "GABRIEL", "AMADO", "MERG", "MARILIA", "ALESSANDRO", "MELISSA", "BALTASAR", "ELISANDRA", "SERGIO", "JOAO", "WAGNER", "MARCIA", "RENATA", "ERICA", "WANDIR", "MARCELO", "FABRICIO", "JUNIO", "RAFAEL", "CARNEIRO", "ANA", "MARCELO", "LEONARDO", "DIEGO", "RUBENS", "EDMUNDO", "MARCO", "FLORDELIZ", "CAROLINA", "RENAN", "CATARINE", "EDNA", "JO", "ZELI", "RICARDO", "ROSANGELA", "SERGIO", "CRISTINA", "HEITOR", "FLASH", "JOAO", "FABRICIO", "JOAO", "BRUNO", "CLISIA", "LUIS", "DANIELI", "RAFAEL", "MARCUS", "RICARDO", "ANA", "MAURO", "JOSE", "VALERIA", "ANDREA", "MONICA", "MATEUS", "ALEXANDRE", "GUILHERME", "VIOMAR", "ELENA", "HEBERT", "CAROLINA", "BRUNA", "CARLOS", "LUISA", "DAIANNE", "JOSE", "ANDREA", "ISABELLA", "GABRIEL", "DAISY", "AGUINALDO", "LUCIANA", "VITOR", "JOAO", "NATALIA", "ISABELLA", "CAROLINA", "ADRIANO", "GABRIEL", "YESHUA", "FRANKLIN", "RAFAEL", "GABRIEL", "ALEXANDRE", "JALISON", "HFR", "BASE", "GUILHERME", "VILMA", "NORTON", "BRUNO")
$values = @(450000, 100000, 60000, 21367.75, 4292, 2609.76, 940.23, 900.21, 752.05, 694.83, 556.35, 524.92, 503.59, 446.18, 364.49, 323.87, 279.96, 191.02, 100.15, 98.96, 97.2, 94.87, 94.24, 93.29, 93.29, 93.1, 92.21, 89.47, 87.61, 86.38, 85.9, 85.59, 83.49, 83.09, 80, 79.87, 79.23, 77.23, 76.01, 74.22, 73.48, 70.94, 70.58, 70.23, 70.02, 68.42, 68.22, 67.76, 67.39, 67.06, 67.03, 66.93, 66.86, 66.47, 64.77, 62.82, 62.24, 61.49, 61.37, 61.18, 61.09, 60.56, 59.71, 59.47, 59.26, 58.71, 57.75, 57.68, 57.48, 57.2, 57.01, 56.88, 56.38, 56.33, 53.8, 53.39, 53.11, 51.44, 51.09, 51.06, 48.95, 48.27, 48.25, 47.9, 46.28, 44.73, 44.17, 43.87, 42.94, 42.61, 42.34, 41.94, 40.34, 40.07, 39.91, 39.09, 37.58, 37.11, 36.75, 36.61, 36.55, 36.46, 35.77, 35.27, 34.71, 33.41, 33.31, 32.34, 31.9, 31.25, 31.01, 30.41, 30, 28.73, 28.18, 28.16, 28.12, 27.83, 27.47, 26.71, 26.65, 25.45, 25.08, 24.96, 24.5, 24.09, 23.36, 22.95, 22.77, 22.74, 22.74, 21.52, 20.89, 20.75, 20.69, 19.2, 18.77, 18.64, 18.19, 17.86, 17.85, 17.62, 17.15, 16.7, 16.48, 16.04, 16.02, 15.62, 15.04, 14.67, 14.59, 14.49, 14.49, 13.38, 13.22, 12.84, 12.67, 12.08, 11.88, 11.8, 11.17, 10.89, 9.79, 9.72, 9.64, 9.12, 8.46, 7.69, 7.54, 7.46, 6.99, 6.71, 6.46, 6.12, 5.86, 5.55, 5.33, 5.32, 5.18, 5.18, 5.17, 4.98, 4.64, 4.4, 4.37, 4.3, 4.2, 3.74, 3.55, 3.11, 2.84, 2.2, 1.7, 1.1, 1, 0.96, 0.88, 0.85, 0.78, 0.69, 0.62, 0.62, 0.5, 0.48, 0.39, 0.29, 0.21, 0.18, 0.17, 0.16, 0.15, 0.14, 0.13, 0.1, 0.09, 0.09, 0.04, 0.04, 0.04, 0.03, 0.02, 0.02, 0.02, 0.01, 0.01, 0.01, 0.01, 0.01)

# Remove the old data rows (row 2 through row 231 -- the 230 original data rows),
# keeping the header row (row 1) and the trailing empty + "Filtros aplicados" rows.
$ws.Rows("2:231").Delete()

# Insert blank rows for the new data set (228 rows) right after the header,
# pushing the trailing empty/filters rows back down.
$lastRow = 1 + $accounts.Length
$ws.Rows("2:" + $lastRow).Insert()

# Make sure the account-number column keeps its leading zeros by formatting it as text
# before writing the values (mirrors the original inline-string storage).
$dataRange = $ws.Range("A2:A" + $lastRow)
$dataRange.NumberFormat = "@"

for ($i = 0; $i -lt $accounts.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $accounts[$i]
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $values[$i]
}
